# New crime data collected - weekly CompStat report refresh
# Updates: report volume/number + week-covering dates, and the weekly
# crime-complaint figures (Week to Date / 28 Day / Year to Date / %Chg
# columns) for precinct rows 15-31 on the active (only) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: volume/number and the "week covering" date range
# ---------------------------------------------------------------------
$ws.Range("A8").Value  = "Volume 31   Number  19"
$ws.Range("C9").Value  = "Report Covering the Week  5/6/2024  Through  5/12/2024"

# ---------------------------------------------------------------------
# Plain value updates (no type change) - row 15 (Rape)
# ---------------------------------------------------------------------
$ws.Range("N15").Value = -77.777777777777

# Row 16 (Robbery)
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = -11.111111111111
$ws.Range("I16").Value = 33
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = -5.714285714285
$ws.Range("L16").Value = -25
$ws.Range("M16").Value = 3.125
$ws.Range("N16").Value = -84.792626728110

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 250
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 41
$ws.Range("J17").Value = 45
$ws.Range("K17").Value = -8.888888888888
$ws.Range("L17").Value = 32.258064516129
$ws.Range("M17").Value = 51.851851851851
$ws.Range("N17").Value = 17.142857142857

# Row 18 (Burglary)
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 63
$ws.Range("J18").Value = 51
$ws.Range("K18").Value = 23.529411764705
$ws.Range("L18").Value = 31.25
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = -77.659574468085

# Row 19 (Gr. Larceny)
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 252
$ws.Range("J19").Value = 255
$ws.Range("K19").Value = -1.176470588235
$ws.Range("L19").Value = -5.617977528089
$ws.Range("M19").Value = 2.857142857142
$ws.Range("N19").Value = -64.556962025316

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 8
$ws.Range("I20").Value = 15
$ws.Range("J20").Value = 30
$ws.Range("L20").Value = -25
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -96.535796766743

# Row 21 (TOTAL)
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -7.407407407407
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = 1.063829787234
$ws.Range("I21").Value = 406
$ws.Range("J21").Value = 422
$ws.Range("K21").Value = -3.791469194312
$ws.Range("L21").Value = -2.168674698795
$ws.Range("M21").Value = 12.777777777777
$ws.Range("N21").Value = -75.947867298578

# ---------------------------------------------------------------------
# Row 22 (Transit) - D22/E22 become "N/A" placeholders (text), rest are
# plain value updates
# ---------------------------------------------------------------------
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "***.*"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("L22").Value = 33.333333333333

# ---------------------------------------------------------------------
# Row 23 (Housing) - D23/E23 become "N/A" placeholders (text), rest are
# plain value updates
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 166.666666666667
$ws.Range("I23").Value = 21
$ws.Range("K23").Value = 16.666666666666
$ws.Range("L23").Value = 425
$ws.Range("M23").Value = 133.333333333333

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 123.076923076923
$ws.Range("F24").Value = 126
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = 41.573033707865
$ws.Range("I24").Value = 432
$ws.Range("J24").Value = 381
$ws.Range("K24").Value = 13.385826771653
$ws.Range("L24").Value = -13.6
$ws.Range("M24").Value = 23.782234957020

# Row 25 (Retail Theft)
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 90
$ws.Range("F25").Value = 92
$ws.Range("G25").Value = 66
$ws.Range("H25").Value = 39.393939393939
$ws.Range("I25").Value = 336
$ws.Range("J25").Value = 303
$ws.Range("K25").Value = 10.891089108910
$ws.Range("L25").Value = -20.567375886524

# Row 26 (Misd. Assault)
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 400
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 62.5
$ws.Range("I26").Value = 83
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 3.75
$ws.Range("L26").Value = 7.792207792207
$ws.Range("M26").Value = -10.752688172043

# ---------------------------------------------------------------------
# Row 27 (UCR Rape*) - D27/E27/G27/H27 flip from "N/A" placeholders to
# real numbers
# ---------------------------------------------------------------------
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = -100
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 1
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = -62.5

# ---------------------------------------------------------------------
# Row 28 (Other Sex Crimes) - C28/D28/E28 flip from real numbers to
# "N/A" placeholders (text)
# ---------------------------------------------------------------------
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "***.*"
$ws.Range("L28").Value = -64.705882352941

# ---------------------------------------------------------------------
# Row 29 (Shooting Vic.) - C29/F29/I29 flip from "N/A" placeholders to
# real numbers
# ---------------------------------------------------------------------
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("C29").Value = 1
$ws.Range("F29").NumberFormat = "#,##0"
$ws.Range("F29").Value = 1
$ws.Range("I29").NumberFormat = "#,##0"
$ws.Range("I29").Value = 1
$ws.Range("K29").Value = -50
$ws.Range("N29").Value = 0

# ---------------------------------------------------------------------
# Row 30 (Shooting Inc.) - same pattern as row 29
# ---------------------------------------------------------------------
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("C30").Value = 1
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("F30").Value = 1
$ws.Range("I30").NumberFormat = "#,##0"
$ws.Range("I30").Value = 1
$ws.Range("K30").Value = -50
$ws.Range("N30").Value = 0

# ---------------------------------------------------------------------
# Row 31 (Hate Crimes) - D31/E31/G31/H31/J31/K31 flip from "N/A"
# placeholders to real numbers; F31/I31/L31 are plain value updates
# ---------------------------------------------------------------------
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("D31").Value = 1
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 3
$ws.Range("G31").NumberFormat = "#,##0"
$ws.Range("G31").Value = 1
$ws.Range("H31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H31").Value = 200
$ws.Range("I31").Value = 4
$ws.Range("J31").NumberFormat = "#,##0"
$ws.Range("J31").Value = 1
$ws.Range("K31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K31").Value = 300
$ws.Range("L31").Value = 100
